$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 895.4286
$ws.Range("I15").Value = 895.4286
$ws.Range("K15").Value = 2686.2858
$ws.Range("M15").Value = -2517.2858
$ws.Range("H53").Value = 156.17647
$ws.Range("I53").Value = 57.5
$ws.Range("J53").Value = 210
$ws.Range("K53").Value = 57.5
$ws.Range("L53").Value = 210
$ws.Range("M53").Value = 579.5
$ws.Range("N53").Value = -1484
$ws.Range("H98").Value = 43019.844
$ws.Range("I98").Value = 1886.3636
$ws.Range("K98").Value = 1886.3636
$ws.Range("M98").Value = -388.3635999999999
$ws.Range("H112").Value = 1564480.5
$ws.Range("I112").Value = 575
$ws.Range("J112").Value = 1787895.6
$ws.Range("K112").Value = 1725
$ws.Range("L112").Value = 5363686.800000001
$ws.Range("M112").Value = -617
$ws.Range("N112").Value = -5365902.800000001
$ws.Range("H122").Value = 43019.844
$ws.Range("I122").Value = 1886.3636
$ws.Range("K122").Value = 5659.0908
$ws.Range("M122").Value = -3209.0908
$ws.Range("H129").Value = 22275.361
$ws.Range("J129").Value = 32234.906
$ws.Range("L129").Value = 96704.71799999999
$ws.Range("N129").Value = -106704.718
$ws.Range("H132").Value = 32161.46
$ws.Range("I132").Value = 20730.723
$ws.Range("J132").Value = 100745.89
$ws.Range("K132").Value = 62192.16900000001
$ws.Range("L132").Value = 302237.67
$ws.Range("M132").Value = -59662.16900000001
$ws.Range("N132").Value = -307297.67
$ws.Range("H137").Value = 1331620.8
$ws.Range("I137").Value = 2409022
$ws.Range("J137").Value = 5588.4614
$ws.Range("K137").Value = 7227066
$ws.Range("L137").Value = 16765.3842
$ws.Range("M137").Value = -7224516
$ws.Range("N137").Value = -21865.3842
$ws.Range("H138").Value = 1918.358
$ws.Range("I138").Value = 1496
$ws.Range("J138").Value = 2129.537
$ws.Range("K138").Value = 4488
$ws.Range("L138").Value = 6388.610999999999
$ws.Range("M138").Value = 652
$ws.Range("N138").Value = -16668.611

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13703.698
$ws.Range("I32").Value = 12547.667
$ws.Range("K32").Value = 12547.667
$ws.Range("M32").Value = -12260.667
$ws.Range("H74").Value = 1902.9773
$ws.Range("I74").Value = 1515.7297
$ws.Range("J74").Value = 3949.8572
$ws.Range("K74").Value = 1515.7297
$ws.Range("L74").Value = 3949.8572
$ws.Range("M74").Value = -641.7297000000001
$ws.Range("N74").Value = -5697.8572
$ws.Range("H77").Value = 1902.9773
$ws.Range("I77").Value = 1515.7297
$ws.Range("J77").Value = 3949.8572
$ws.Range("K77").Value = 7578.6485
$ws.Range("L77").Value = 19749.286
$ws.Range("M77").Value = -3210.6485
$ws.Range("N77").Value = -28485.286
$ws.Range("H80").Value = 49500
$ws.Range("I80").Value = 10000
$ws.Range("J80").Value = 89000
$ws.Range("K80").Value = 10000
$ws.Range("L80").Value = 89000
$ws.Range("M80").Value = -9002
$ws.Range("N80").Value = -90996
$ws.Range("H83").Value = 49500
$ws.Range("I83").Value = 10000
$ws.Range("J83").Value = 89000
$ws.Range("K83").Value = 30000
$ws.Range("L83").Value = 267000
$ws.Range("M83").Value = -25008
$ws.Range("N83").Value = -276984
$ws.Range("H122").Value = 2011.8334
$ws.Range("I122").Value = 1752.3334
$ws.Range("J122").Value = 2271.3333
$ws.Range("K122").Value = 5257.0002
$ws.Range("L122").Value = 6813.999899999999
$ws.Range("M122").Value = -2807.0002
$ws.Range("N122").Value = -11713.9999
$ws.Range("H132").Value = 20835662
$ws.Range("I132").Value = 35716200
$ws.Range("J132").Value = 2908.4
$ws.Range("K132").Value = 107148600
$ws.Range("L132").Value = 8725.200000000001
$ws.Range("M132").Value = -107146070
$ws.Range("N132").Value = -13785.2

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2237.7886
$ws.Range("I134").Value = 1764.7
$ws.Range("J134").Value = 3814.75
$ws.Range("K134").Value = 5294.1
$ws.Range("L134").Value = 11444.25
$ws.Range("M134").Value = -2759.1
$ws.Range("N134").Value = -16514.25

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3985.7144
$ws.Range("J6").Value = 4225
$ws.Range("L6").Value = 4225
$ws.Range("N6").Value = -4451
$ws.Range("H10").Value = 125300.625
$ws.Range("I10").Value = 200277.8
$ws.Range("J10").Value = 338.66666
$ws.Range("K10").Value = 200277.8
$ws.Range("L10").Value = 338.66666
$ws.Range("M10").Value = -200138.8
$ws.Range("N10").Value = -616.66666
$ws.Range("H132").Value = 1756046.1
$ws.Range("I132").Value = 2980
$ws.Range("K132").Value = 8940
$ws.Range("M132").Value = -6410
$ws.Range("H137").Value = 59663
$ws.Range("J137").Value = 59663
$ws.Range("L137").Value = 59663
$ws.Range("N137").Value = -69863

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1230.9584
$ws.Range("J68").Value = 1319.9814
$ws.Range("L68").Value = 3959.9442
$ws.Range("N68").Value = -5581.9442
$ws.Range("H71").Value = 1230.9584
$ws.Range("J71").Value = 1319.9814
$ws.Range("L71").Value = 11879.8326
$ws.Range("N71").Value = -19991.8326
$ws.Range("H107").Value = 691.2121
$ws.Range("I107").Value = 331.73914
$ws.Range("K107").Value = 995.2174200000001
$ws.Range("M107").Value = 924.7825799999999
$ws.Range("H113").Value = 4767.32
$ws.Range("I113").Value = 8906.583000000001
$ws.Range("K113").Value = 26719.749
$ws.Range("M113").Value = -24549.749
$ws.Range("H122").Value = 3029.9333
$ws.Range("I122").Value = 660.71875
$ws.Range("J122").Value = 8861.846
$ws.Range("K122").Value = 5946.46875
$ws.Range("L122").Value = 79756.614
$ws.Range("M122").Value = -3496.46875
$ws.Range("N122").Value = -84656.614

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 40700
$ws.Range("J74").Value = 40700
$ws.Range("L74").Value = 40700
$ws.Range("N74").Value = -42572
$ws.Range("H77").Value = 40700
$ws.Range("J77").Value = 40700
$ws.Range("L77").Value = 122100
$ws.Range("N77").Value = -131460
$ws.Range("H113").Value = 1977.9286
$ws.Range("I113").Value = 1931.1
$ws.Range("J113").Value = 2095
$ws.Range("K113").Value = 1931.1
$ws.Range("L113").Value = 2095
$ws.Range("M113").Value = 238.9000000000001
$ws.Range("N113").Value = -6435
$ws.Range("H122").Value = 2666
$ws.Range("J122").Value = 2666
$ws.Range("L122").Value = 7998
$ws.Range("N122").Value = -12898
$ws.Range("H126").Value = 9454.6
$ws.Range("I126").Value = 15219.5
$ws.Range("J126").Value = 2866.1428
$ws.Range("K126").Value = 45658.5
$ws.Range("L126").Value = 8598.428400000001
$ws.Range("M126").Value = -43188.5
$ws.Range("N126").Value = -13538.4284
$ws.Range("H132").Value = 58827240
$ws.Range("I132").Value = 142859490
$ws.Range("J132").Value = 4670.5
$ws.Range("K132").Value = 428578470
$ws.Range("L132").Value = 14011.5
$ws.Range("M132").Value = -428575940
$ws.Range("N132").Value = -19071.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 5778.1113
$ws.Range("I40").Value = 4250.25
$ws.Range("J40").Value = 7000.4
$ws.Range("K40").Value = 4250.25
$ws.Range("L40").Value = 7000.4
$ws.Range("M40").Value = -4114.25
$ws.Range("N40").Value = -7272.4
$ws.Range("H55").Value = 505.69565
$ws.Range("I55").Value = 416.66666
$ws.Range("J55").Value = 602.8182
$ws.Range("K55").Value = 416.66666
$ws.Range("L55").Value = 602.8182
$ws.Range("M55").Value = -243.66666
$ws.Range("N55").Value = -948.8182
$ws.Range("H132").Value = 5671.2144
$ws.Range("I132").Value = 5880
$ws.Range("J132").Value = 5555.222
$ws.Range("K132").Value = 17640
$ws.Range("L132").Value = 16665.666
$ws.Range("M132").Value = -15110
$ws.Range("N132").Value = -21725.666
$ws.Range("H136").Value = 2085.3438
$ws.Range("I136").Value = 1585.84
$ws.Range("J136").Value = 3869.2856
$ws.Range("K136").Value = 4757.52
$ws.Range("L136").Value = 11607.8568
$ws.Range("M136").Value = -2207.52
$ws.Range("N136").Value = -16707.8568

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 971.95
$ws.Range("J113").Value = 1422.5
$ws.Range("L113").Value = 4267.5
$ws.Range("N113").Value = -8607.5
$ws.Range("H122").Value = 28571428
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H123").Value = 35200
$ws.Range("J123").Value = 35200
$ws.Range("L123").Value = 35200
$ws.Range("N123").Value = -45000
$ws.Range("H126").Value = 1839842.1
$ws.Range("I126").Value = 2675278
$ws.Range("J126").Value = 1883
$ws.Range("K126").Value = 8025834
$ws.Range("L126").Value = 5649
$ws.Range("M126").Value = -8023364
$ws.Range("N126").Value = -10589
$ws.Range("H127").Value = 39945.453
$ws.Range("J127").Value = 39945.453
$ws.Range("L127").Value = 39945.453
$ws.Range("N127").Value = -49865.453
$ws.Range("H132").Value = 1451406.2
$ws.Range("I132").Value = 2718893.5
$ws.Range("J132").Value = 2849.3572
$ws.Range("K132").Value = 8156680.5
$ws.Range("L132").Value = 8548.071599999999
$ws.Range("M132").Value = -8154150.5
$ws.Range("N132").Value = -13608.0716
$ws.Range("H136").Value = 730360.5600000001
$ws.Range("I136").Value = 1015467.8
$ws.Range("J136").Value = 1753.1111
$ws.Range("K136").Value = 3046403.4
$ws.Range("L136").Value = 5259.3333
$ws.Range("M136").Value = -3043853.4
$ws.Range("N136").Value = -10359.3333
$ws.Range("N122").ClearContents()
